$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 576, shifting existing rows (576-617) down to (577-618)
$ws.Rows.Item(576).Insert()

# Populate the newly inserted row 576 with the new data point.
# Leading apostrophe forces the date-looking text to stay a literal string
# (matches the rest of column A, which stores dates as plain text), then
# restore the default "Normal" style so no stray quote-prefix formatting
# is left behind on the cell.
$ws.Cells.Item(576, 1).Value = "'2026/01/06"
$ws.Cells.Item(576, 1).Style = "Normal"
$ws.Cells.Item(576, 2).Value = "火"
$ws.Cells.Item(576, 3).Value = 12
$ws.Cells.Item(576, 4).Value = 201
